$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.049133904040056
$ws.Range("D2").Value = 1.057344047581932
$ws.Range("E2").Value = 1.052807067124932
$ws.Range("F2").Value = 1.066030348053597
$ws.Range("I2").Value = 1.050574343273791
$ws.Range("J2").Value = 1.054173963831215
$ws.Range("K2").Value = 1.060079274988282
$ws.Range("L2").Value = 1.055554779218594
$ws.Range("M2").Value = 1.068741993696972
$ws.Range("N2").Value = 1.055671011775299

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.050128437672494
$ws.Range("D3").Value = 1.05815795285234
$ws.Range("E3").Value = 1.053751092747251
$ws.Range("F3").Value = 1.066975360151392
$ws.Range("I3").Value = 1.05090203186054
$ws.Range("J3").Value = 1.054817434384725
$ws.Range("K3").Value = 1.060706997738683
$ws.Range("L3").Value = 1.056311387722772
$ws.Range("M3").Value = 1.069502195958205
$ws.Range("N3").Value = 1.056315396130802

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.050772315414106
$ws.Range("D4").Value = 1.058684881767409
$ws.Range("E4").Value = 1.054362614179068
$ws.Range("F4").Value = 1.067587448745815
$ws.Range("I4").Value = 1.051112974876442
$ws.Range("J4").Value = 1.055233530099396
$ws.Range("K4").Value = 1.061112800653137
$ws.Range("L4").Value = 1.056801011761477
$ws.Range("M4").Value = 1.069994071778729
$ws.Range("N4").Value = 1.056732082749115

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.051043084364872
$ws.Range("D5").Value = 1.058906468186454
$ws.Range("E5").Value = 1.054619858150673
$ws.Range("F5").Value = 1.067844914039787
$ws.Range("I5").Value = 1.0512013931289
$ws.Range("J5").Value = 1.055408390773451
$ws.Range("K5").Value = 1.061283309536466
$ws.Range("L5").Value = 1.057006860368716
$ws.Range("M5").Value = 1.070200849338638
$ws.Range("N5").Value = 1.056907191745384

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.051088552499424
$ws.Range("D6").Value = 1.058943677333391
$ws.Range("E6").Value = 1.054663059922136
$ws.Range("F6").Value = 1.067888151966453
$ws.Range("I6").Value = 1.051216223555175
$ws.Range("J6").Value = 1.055437746762595
$ws.Range("K6").Value = 1.061311933390813
$ws.Range("L6").Value = 1.057041423858555
$ws.Range("M6").Value = 1.070235567763899
$ws.Range("N6").Value = 1.0569365894234

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.050775933117619
$ws.Range("D7").Value = 1.05868784235985
$ws.Range("E7").Value = 1.054366050856654
$ws.Range("F7").Value = 1.067590888448062
$ws.Range("I7").Value = 1.051114157356219
$ws.Range("J7").Value = 1.055235866856303
$ws.Range("K7").Value = 1.061115079358322
$ws.Range("L7").Value = 1.056803762279706
$ws.Range("M7").Value = 1.069996834779296
$ws.Range("N7").Value = 1.056734422824485

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.049469938949012
$ws.Range("D8").Value = 1.0576190525621
$ws.Range("E8").Value = 1.053125964943895
$ws.Range("F8").Value = 1.066349593942094
$ws.Range("I8").Value = 1.050685313306966
$ws.Range("J8").Value = 1.054391483794214
$ws.Range("K8").Value = 1.060291494182326
$ws.Range("L8").Value = 1.055810468256769
$ws.Range("M8").Value = 1.068998912462439
$ws.Range("N8").Value = 1.055888840641595

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.047171300013344
$ws.Range("D9").Value = 1.055737879813348
$ws.Range("E9").Value = 1.050945976066537
$ws.Range("F9").Value = 1.0641669338909
$ws.Range("I9").Value = 1.04992127987494
$ws.Range("J9").Value = 1.052901513198555
$ws.Range("K9").Value = 1.058837386472172
$ws.Range("L9").Value = 1.054060556778149
$ws.Range("M9").Value = 1.067240284671155
$ws.Range("N9").Value = 1.054396754116788

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.045640713572921
$ws.Range("D10").Value = 1.054485282662931
$ws.Range("E10").Value = 1.0494962007988
$ws.Range("F10").Value = 1.062715018492525
$ws.Range("I10").Value = 1.049406331099489
$ws.Range("J10").Value = 1.051906852876932
$ws.Range("K10").Value = 1.057866110054144
$ws.Range("L10").Value = 1.052894262066795
$ws.Range("M10").Value = 1.066067806419356
$ws.Range("N10").Value = 1.05340068126342

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.044978394978952
$ws.Range("D11").Value = 1.05394326617578
$ws.Range("E11").Value = 1.048869285158477
$ws.Range("F11").Value = 1.062087091924426
$ws.Range("I11").Value = 1.049182030877898
$ws.Range("J11").Value = 1.051475842258801
$ws.Range("K11").Value = 1.057445101238188
$ws.Range("L11").Value = 1.05238932673279
$ws.Range("M11").Value = 1.065560106066435
$ws.Range("N11").Value = 1.052969058560773

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04473244617875
$ws.Range("D12").Value = 1.053741993065954
$ws.Range("E12").Value = 1.04863654879489
$ws.Range("F12").Value = 1.061853967483408
$ws.Range("I12").Value = 1.04909851721868
$ws.Range("J12").Value = 1.051315698792578
$ws.Range("K12").Value = 1.057288654382344
$ws.Range("L12").Value = 1.0522017836462
$ws.Range("M12").Value = 1.065371522855955
$ws.Range("N12").Value = 1.05280868767246

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.044785200046588
$ws.Range("D13").Value = 1.053785164284957
$ws.Range("E13").Value = 1.048686465728555
$ws.Range("F13").Value = 1.061903968228682
$ws.Range("I13").Value = 1.049116440162188
$ws.Range("J13").Value = 1.051350052224391
$ws.Range("K13").Value = 1.057322215715039
$ws.Range("L13").Value = 1.05224201170405
$ws.Range("M13").Value = 1.065411974614296
$ws.Range("N13").Value = 1.052843089890087

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.044958063429478
$ws.Range("D14").Value = 1.053926627731183
$ws.Range("E14").Value = 1.048850044481745
$ws.Range("F14").Value = 1.062067819425122
$ws.Range("I14").Value = 1.049175131661525
$ws.Range("J14").Value = 1.051462605704051
$ws.Range("K14").Value = 1.057432170623526
$ws.Range("L14").Value = 1.052373824106374
$ws.Range("M14").Value = 1.065544517729308
$ws.Range("N14").Value = 1.052955803208597

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.045064578998011
$ws.Range("D15").Value = 1.054013795469949
$ws.Range("E15").Value = 1.048950847738769
$ws.Range("F15").Value = 1.062168788879712
$ws.Range("I15").Value = 1.049211267129309
$ws.Range("J15").Value = 1.051531947415585
$ws.Range("K15").Value = 1.057499908819467
$ws.Range("L15").Value = 1.05245503973802
$ws.Range("M15").Value = 1.065626181847098
$ws.Range("N15").Value = 1.053025243393316

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.04568467868451
$ws.Range("D16").Value = 1.054521262308007
$ws.Range("E16").Value = 1.049537825032683
$ws.Range("F16").Value = 1.062756708050569
$ws.Range("I16").Value = 1.049421189300763
$ws.Range("J16").Value = 1.051935451019621
$ws.Range("K16").Value = 1.057894041811773
$ws.Range("L16").Value = 1.052927774661375
$ws.Range("M16").Value = 1.066101500681817
$ws.Range("N16").Value = 1.053429320018752

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.046073767691398
$ws.Range("D17").Value = 1.0548396815698
$ws.Range("E17").Value = 1.049906247625958
$ws.Range("F17").Value = 1.063125698963575
$ws.Range("I17").Value = 1.04955251370327
$ws.Range("J17").Value = 1.052188473928955
$ws.Range("K17").Value = 1.0581411537417
$ws.Range("L17").Value = 1.05322433013291
$ws.Range("M17").Value = 1.06639965359733
$ws.Range("N17").Value = 1.053682702249639

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.046300758712568
$ws.Range("D18").Value = 1.055025445424029
$ws.Range("E18").Value = 1.050121223922351
$ws.Range("F18").Value = 1.063340998618127
$ws.Range("I18").Value = 1.049628985234852
$ws.Range("J18").Value = 1.052336027370757
$ws.Range("K18").Value = 1.058285247481134
$ws.Range("L18").Value = 1.05339731340066
$ws.Range("M18").Value = 1.066573560179211
$ws.Range("N18").Value = 1.053830465234253

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.046378163894446
$ws.Range("D19").Value = 1.055088792055279
$ws.Range("E19").Value = 1.050194539104928
$ws.Range("F19").Value = 1.063414422675464
$ws.Range("I19").Value = 1.049655038368315
$ws.Range("J19").Value = 1.05238633408399
$ws.Range("K19").Value = 1.058334372500147
$ws.Range("L19").Value = 1.053456297472989
$ws.Range("M19").Value = 1.066632857642633
$ws.Range("N19").Value = 1.053880843388789

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.046032017777487
$ws.Range("D20").Value = 1.0548055145401
$ws.Range("E20").Value = 1.049866710918942
$ws.Range("F20").Value = 1.063086102121889
$ws.Range("I20").Value = 1.049538437051101
$ws.Range("J20").Value = 1.05216133012019
$ws.Range("K20").Value = 1.058114645358412
$ws.Range("L20").Value = 1.053192511759004
$ws.Range("M20").Value = 1.066367664706078
$ws.Range("N20").Value = 1.053655519893551

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.044907157651689
$ws.Range("D21").Value = 1.053884968720454
$ws.Range("E21").Value = 1.048801871075453
$ws.Range("F21").Value = 1.062019566137529
$ws.Range("I21").Value = 1.049157853956521
$ws.Range("J21").Value = 1.051429462800623
$ws.Range("K21").Value = 1.057399793444848
$ws.Range("L21").Value = 1.052335008292168
$ws.Range("M21").Value = 1.065505487098708
$ws.Range("N21").Value = 1.052922613238445

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.044200294432551
$ws.Range("D22").Value = 1.053306509168306
$ws.Range("E22").Value = 1.048133105088649
$ws.Range("F22").Value = 1.061349661011519
$ws.Range("I22").Value = 1.048917417569021
$ws.Range("J22").Value = 1.050969037516672
$ws.Range("K22").Value = 1.056949959392661
$ws.Range("L22").Value = 1.051795933778517
$ws.Range("M22").Value = 1.064963398168868
$ws.Range("N22").Value = 1.052461534097779

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.044574979752355
$ws.Range("D23").Value = 1.053613130493539
$ws.Range("E23").Value = 1.04848756006955
$ws.Range("F23").Value = 1.061704726736322
$ws.Range("I23").Value = 1.0490449862133
$ws.Range("J23").Value = 1.05121314319387
$ws.Range("K23").Value = 1.057188460581846
$ws.Range("L23").Value = 1.052081700335308
$ws.Range("M23").Value = 1.065250769865971
$ws.Range("N23").Value = 1.052705986433038

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.046050882636699
$ws.Range("D24").Value = 1.054820953041144
$ws.Range("E24").Value = 1.049884575602811
$ws.Range("F24").Value = 1.063103994003169
$ws.Range("I24").Value = 1.04954479807908
$ws.Range("J24").Value = 1.05217359533335
$ws.Range("K24").Value = 1.058126623487017
$ws.Range("L24").Value = 1.053206889089474
$ws.Range("M24").Value = 1.06638211911211
$ws.Range("N24").Value = 1.053667802524721

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.047765231347293
$ws.Range("D25").Value = 1.056223944754047
$ws.Range("E25").Value = 1.0515089335679
$ws.Range("F25").Value = 1.064730645635977
$ws.Range("I25").Value = 1.050119788470614
$ws.Range("J25").Value = 1.053286946467964
$ws.Range("K25").Value = 1.059213641063502
$ws.Range("L25").Value = 1.054512898936861
$ws.Range("M25").Value = 1.067694946204032
$ws.Range("N25").Value = 1.05478273474565
